$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

# Row 21
$ws.Range("H21").Value = 21500
$ws.Range("I21").Value = 18333.334
$ws.Range("K21").Value = 18333.334
$ws.Range("M21").Value = -17865.334

# Row 23
$ws.Range("H23").Value = 21500
$ws.Range("I23").Value = 18333.334
$ws.Range("K23").Value = 18333.334
$ws.Range("M23").Value = -18099.334

# Row 29
$ws.Range("H29").Value = 892.8570999999999

# Row 43
$ws.Range("H43").Value = 575
$ws.Range("J43").Value = 575
$ws.Range("L43").Value = 575
$ws.Range("N43").Value = -713

# Row 92
$ws.Range("H92").Value = 2147.8
$ws.Range("I92").Value = 2515.9
$ws.Range("J92").Value = 675.4
$ws.Range("K92").Value = 2515.9
$ws.Range("L92").Value = 675.4
$ws.Range("M92").Value = -1267.9
$ws.Range("N92").Value = -3171.4

# Row 98
$ws.Range("H98").Value = 293.9
$ws.Range("I98").Value = 293.9
$ws.Range("K98").Value = 293.9
$ws.Range("M98").Value = 1204.1

# Row 100
$ws.Range("H100").Value = 2352.1428
$ws.Range("I100").Value = 1716.25
$ws.Range("J100").Value = 3200
$ws.Range("K100").Value = 1716.25
$ws.Range("L100").Value = 3200
$ws.Range("M100").Value = -1175.25
$ws.Range("N100").Value = -4282

# Row 122
$ws.Range("H122").Value = 293.9
$ws.Range("I122").Value = 293.9
$ws.Range("K122").Value = 881.6999999999999
$ws.Range("M122").Value = 1568.3

# Row 132
$ws.Range("H132").Value = 35934.934
$ws.Range("I132").Value = 37160.277
$ws.Range("K132").Value = 111480.831
$ws.Range("M132").Value = -108950.831

# Row 135
$ws.Range("H135").Value = 4284.8237
$ws.Range("I135").Value = 1238.6666
$ws.Range("J135").Value = 11595.6
$ws.Range("K135").Value = 11147.9994
$ws.Range("L135").Value = 104360.4
$ws.Range("M135").Value = -8612.999400000001
$ws.Range("N135").Value = -109430.4

# Row 137
$ws.Range("H137").Value = 2176.7097
$ws.Range("I137").Value = 2469.6667
$ws.Range("K137").Value = 7409.000100000001
$ws.Range("M137").Value = -4859.000100000001

# Row 138
$ws.Range("H138").Value = 2476.2068
$ws.Range("I138").Value = 704.9167
$ws.Range("J138").Value = 3726.5293
$ws.Range("K138").Value = 2114.7501
$ws.Range("L138").Value = 11179.5879
$ws.Range("M138").Value = 3025.2499
$ws.Range("N138").Value = -21459.5879

# Row 140
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = $null

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

# Row 2
$ws.Range("H2").Value = 3994.7144
$ws.Range("I2").Value = 2475
$ws.Range("J2").Value = 4602.6
$ws.Range("K2").Value = 2475
$ws.Range("L2").Value = 4602.6
$ws.Range("M2").Value = -2362
$ws.Range("N2").Value = -4828.6

# Row 32
$ws.Range("H32").Value = 2308.2188
$ws.Range("I32").Value = 1400.2549
$ws.Range("K32").Value = 1400.2549
$ws.Range("M32").Value = -1113.2549

# Row 61
$ws.Range("H61").Value = 3292
$ws.Range("I61").Value = 2777.7
$ws.Range("J61").Value = 4434.8887
$ws.Range("K61").Value = 2777.7
$ws.Range("L61").Value = 4434.8887
$ws.Range("M61").Value = -2565.7
$ws.Range("N61").Value = -4858.8887

# Row 74
$ws.Range("H74").Value = 815.6070999999999
$ws.Range("I74").Value = 332.7647
$ws.Range("J74").Value = 1561.8182
$ws.Range("K74").Value = 332.7647
$ws.Range("L74").Value = 1561.8182
$ws.Range("M74").Value = 541.2353000000001
$ws.Range("N74").Value = -3309.8182

# Row 77
$ws.Range("H77").Value = 815.6070999999999
$ws.Range("I77").Value = 332.7647
$ws.Range("J77").Value = 1561.8182
$ws.Range("K77").Value = 1663.8235
$ws.Range("L77").Value = 7809.090999999999
$ws.Range("M77").Value = 2704.1765
$ws.Range("N77").Value = -16545.091

# Row 116
$ws.Range("H116").Value = 3994.7144
$ws.Range("I116").Value = 2475
$ws.Range("J116").Value = 4602.6
$ws.Range("K116").Value = 2475
$ws.Range("L116").Value = 4602.6
$ws.Range("M116").Value = -181
$ws.Range("N116").Value = -9190.6

# Row 136
$ws.Range("H136").Value = 3292
$ws.Range("I136").Value = 2777.7
$ws.Range("J136").Value = 4434.8887
$ws.Range("K136").Value = 8333.099999999999
$ws.Range("L136").Value = 13304.6661
$ws.Range("M136").Value = -5783.099999999999
$ws.Range("N136").Value = -18404.6661

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

# Row 3
$ws.Range("H3").Value = 3994.7144
$ws.Range("I3").Value = 2475
$ws.Range("J3").Value = 4602.6
$ws.Range("K3").Value = 2475
$ws.Range("L3").Value = 4602.6
$ws.Range("M3").Value = -2361
$ws.Range("N3").Value = -4830.6

# Row 81
$ws.Range("H81").Value = 15985.5
$ws.Range("J81").Value = 15985.5
$ws.Range("L81").Value = 15985.5
$ws.Range("N81").Value = -18107.5

# Row 84
$ws.Range("H84").Value = 15985.5
$ws.Range("J84").Value = 15985.5
$ws.Range("L84").Value = 47956.5
$ws.Range("N84").Value = -58564.5

# Row 86
$ws.Range("H86").Value = 1657.8
$ws.Range("I86").Value = 1325
$ws.Range("J86").Value = 2157
$ws.Range("K86").Value = 1325
$ws.Range("L86").Value = 2157
$ws.Range("M86").Value = -202
$ws.Range("N86").Value = -4403

# Row 89
$ws.Range("H89").Value = 1657.8
$ws.Range("I89").Value = 1325
$ws.Range("J89").Value = 2157
$ws.Range("K89").Value = 6625
$ws.Range("L89").Value = 10785
$ws.Range("M89").Value = -1009
$ws.Range("N89").Value = -22017

# Row 99
$ws.Range("H99").Value = 830
$ws.Range("I99").Value = 816.1875
$ws.Range("K99").Value = 816.1875
$ws.Range("M99").Value = 681.8125

# Row 107
$ws.Range("H107").Value = 1184.25
$ws.Range("I107").Value = 1023.44446
$ws.Range("K107").Value = 1023.44446
$ws.Range("M107").Value = 896.55554

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

# Row 31
$ws.Range("H31").Value = 7845.0835
$ws.Range("I31").Value = 10727.4
$ws.Range("J31").Value = 3041.2222
$ws.Range("K31").Value = 10727.4
$ws.Range("L31").Value = 3041.2222
$ws.Range("M31").Value = -10432.4
$ws.Range("N31").Value = -3631.2222

# Row 34
$ws.Range("H34").Value = 7845.0835
$ws.Range("I34").Value = 10727.4
$ws.Range("J34").Value = 3041.2222
$ws.Range("K34").Value = 10727.4
$ws.Range("L34").Value = 3041.2222
$ws.Range("M34").Value = -10525.4
$ws.Range("N34").Value = -3445.2222

# Row 38
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").Value = $null

# Row 46
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").Value = $null

# Row 107
$ws.Range("H107").Value = 353.625
$ws.Range("I107").Value = 332.7143
$ws.Range("J107").Value = 500
$ws.Range("K107").Value = 332.7143
$ws.Range("L107").Value = 500
$ws.Range("M107").Value = 1587.2857
$ws.Range("N107").Value = -4340

# Row 132
$ws.Range("H132").Value = 18687.773
$ws.Range("I132").Value = 25884.85
$ws.Range("J132").Value = 5602.1816
$ws.Range("K132").Value = 77654.54999999999
$ws.Range("L132").Value = 16806.5448
$ws.Range("M132").Value = -75124.54999999999
$ws.Range("N132").Value = -21866.5448

# Row 134
$ws.Range("H134").Value = 790.63635
$ws.Range("I134").Value = 730.2353000000001
$ws.Range("J134").Value = 996
$ws.Range("K134").Value = 2190.7059
$ws.Range("L134").Value = 2988
$ws.Range("M134").Value = 344.2941000000001
$ws.Range("N134").Value = -8058

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

# Row 5
$ws.Range("H5").Value = 770.2917
$ws.Range("I5").Value = 701.3333
$ws.Range("J5").Value = 780.1429000000001
$ws.Range("K5").Value = 2103.9999
$ws.Range("L5").Value = 2340.4287
$ws.Range("M5").Value = -1991.9999
$ws.Range("N5").Value = -2564.4287

# Row 33
$ws.Range("H33").Value = 421
$ws.Range("I33").Value = 452.5
$ws.Range("J33").Value = 295
$ws.Range("K33").Value = 2715
$ws.Range("L33").Value = 1770
$ws.Range("M33").Value = -2432
$ws.Range("N33").Value = -2336

# Row 131
$ws.Range("H131").Value = 805.25
$ws.Range("I131").Value = 405.66666
$ws.Range("J131").Value = 830.7553
$ws.Range("K131").Value = 1216.99998
$ws.Range("L131").Value = 2492.2659
$ws.Range("M131").Value = 3823.00002
$ws.Range("N131").Value = -12572.2659

# Row 132
$ws.Range("H132").Value = 1343.1428
$ws.Range("I132").Value = 600.5
$ws.Range("J132").Value = 2333.3333
$ws.Range("K132").Value = 5404.5
$ws.Range("L132").Value = 20999.9997
$ws.Range("M132").Value = -2874.5
$ws.Range("N132").Value = -26059.9997

# Row 135
$ws.Range("H135").Value = 770.2917
$ws.Range("I135").Value = 701.3333
$ws.Range("J135").Value = 780.1429000000001
$ws.Range("K135").Value = 6311.9997
$ws.Range("L135").Value = 7021.2861
$ws.Range("M135").Value = -3776.9997
$ws.Range("N135").Value = -12091.2861

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

# Row 70
$ws.Range("H70").Value = 4842
$ws.Range("I70").Value = 4300
$ws.Range("J70").Value = 4977.5
$ws.Range("K70").Value = 4300
$ws.Range("L70").Value = 4977.5
$ws.Range("M70").Value = -4030
$ws.Range("N70").Value = -5517.5

# Row 73
$ws.Range("H73").Value = 4842
$ws.Range("I73").Value = 4300
$ws.Range("J73").Value = 4977.5
$ws.Range("K73").Value = 4300
$ws.Range("L73").Value = 4977.5
$ws.Range("M73").Value = -3364
$ws.Range("N73").Value = -6849.5

# Row 97
$ws.Range("H97").Value = 1784.68
$ws.Range("I97").Value = 1012.1177
$ws.Range("J97").Value = 3426.375
$ws.Range("K97").Value = 1012.1177
$ws.Range("L97").Value = 3426.375
$ws.Range("M97").Value = -516.1177
$ws.Range("N97").Value = -4418.375

# Row 107
$ws.Range("H107").Value = 1494.6
$ws.Range("I107").Value = 391
$ws.Range("J107").Value = 3150
$ws.Range("K107").Value = 391
$ws.Range("L107").Value = 3150
$ws.Range("M107").Value = 1529
$ws.Range("N107").Value = -6990

# Row 122
$ws.Range("H122").Value = 2628.4546
$ws.Range("I122").Value = 2254.75
$ws.Range("J122").Value = 3625
$ws.Range("K122").Value = 6764.25
$ws.Range("L122").Value = 10875
$ws.Range("M122").Value = -4314.25
$ws.Range("N122").Value = -15775

# Row 132
$ws.Range("H132").Value = 22766.375
$ws.Range("I132").Value = 1366.3889
$ws.Range("J132").Value = 86966.336
$ws.Range("K132").Value = 4099.1667
$ws.Range("L132").Value = 260899.008
$ws.Range("M132").Value = -1569.1667
$ws.Range("N132").Value = -265959.008

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

# Row 34
$ws.Range("H34").Value = 19937.166
$ws.Range("I34").Value = 9919.799999999999
$ws.Range("J34").Value = 70024
$ws.Range("K34").Value = 9919.799999999999
$ws.Range("L34").Value = 70024
$ws.Range("M34").Value = -9747.799999999999
$ws.Range("N34").Value = -70368

# Row 61
$ws.Range("H61").Value = 8145.75
$ws.Range("I61").Value = 4548.8
$ws.Range("J61").Value = 10715
$ws.Range("K61").Value = 4548.8
$ws.Range("L61").Value = 10715
$ws.Range("M61").Value = -4346.8
$ws.Range("N61").Value = -11119

# Row 113
$ws.Range("H113").Value = 8145.75
$ws.Range("I113").Value = 4548.8
$ws.Range("J113").Value = 10715
$ws.Range("K113").Value = 4548.8
$ws.Range("L113").Value = 10715
$ws.Range("M113").Value = -2378.8
$ws.Range("N113").Value = -15055
